$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.555.76"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "2.383.31"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("D5").Value = "'507.10"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'131.21"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").Value = "2.396.86"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("E10").Value = "  +2.54%  "

$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("D12").Value = "'4.84"
$ws.Range("E12").Value = "  +5.87%  "

$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").Value = "2.806.63"
$ws.Range("E14").Value = "  -0.89%  "

$ws.Range("D15").Value = "56.515.99"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "2.453.17"
$ws.Range("E18").Value = "  +3.85%  "

$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "'310.53"
$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D24").Value = "'66.53"
$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("E27").Value = "  -3.44%  "

$ws.Range("D28").Value = "'7.28"
$ws.Range("E28").Value = "  -2.65%  "

$ws.Range("D29").Value = "'173.36"
$ws.Range("E29").Value = "  +1.33%  "

$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.08"
$ws.Range("E34").Value = "  -3.57%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.995"
$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("E36").Value = "  -1.10%  "

$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("E38").Value = "  -2.39%  "

$ws.Range("D39").Value = "'0.832"
$ws.Range("E39").Value = "  +3.77%  "

$ws.Range("D40").Value = "'36.59"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("E41").Value = "  -2.89%  "

$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "  +1.45%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'4.97"
$ws.Range("E43").Value = "  +1.94%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'129.08"
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("E46").Value = "  -0.94%  "

$ws.Range("E47").Value = "  -4.17%  "

$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("E49").Value = "  -0.64%  "

$ws.Range("D50").Value = "'17.18"
$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("E51").Value = "  -1.05%  "

